# edit.ps1 - apply the two textual changes described by the diff:
#  1. Slide 11, body placeholder ("Rectangle 3"), paragraph 9:
#     "can also run batch files" -> "can also run command files (a.k.a. batch files)"
#  2. Slide 7, body placeholder ("Rectangle 3"), paragraph 6:
#     "Testing/Re-engineering tools" -> "Testing/reverse-engineering tools"
#     (re-split across 4 runs: "Testing/" | "r" | "everse-engineering " | "tools")

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Change 1 - Slide 11 ("Examples of Interpreters")
# ---------------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$shape11 = $slide11.Shapes.Item(4)          # "Rectangle 3" body placeholder
$tr11 = $shape11.TextFrame.TextRange
$para11 = $tr11.Paragraphs(9, 1)            # "can also run batch files"

# Route the text change through an intermediate value that shares no
# substring with either the old or new text, so the engine's prefix/suffix
# "diff merge" doesn't split the run in two - the final assignment then
# lands cleanly as a single run with the original run properties intact.
$para11.Text = "QzQzPlaceholderQzQz"
$para11.Text = "can also run command files (a.k.a. batch files)"

# ---------------------------------------------------------------------------
# Change 2 - Slide 7 ("Other Language Processors")
# ---------------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$shape7 = $slide7.Shapes.Item(2)            # "Rectangle 3" body placeholder
$tr7 = $shape7.TextFrame.TextRange
$para7 = $tr7.Paragraphs(6, 1)              # "Testing/Re-engineering tools"

# Same trick to land the whole new sentence as a single run first.
$para7.Text = "QzQzPlaceholderQzQz"
$para7.Text = "Testing/reverse-engineering tools"

# Now split that single run into the four runs shown in the diff:
#   "Testing/" + "r" + "everse-engineering " + "tools"
$run1 = $para7.Characters(1, 8)
$run1.Text = "Testing/"

$run2 = $para7.Characters(9, 1)
$run2.Text = "r"

$run3 = $para7.Characters(10, 19)
$run3.Text = "everse-engineering "

$run4 = $para7.Characters(29, 5)
$run4.Text = "tools"
